$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Reliance"
$ws.Range("B1").Value = "Niva Bupa (formerly known as Max Bupa)"
$ws.Range("C1").Value = "Care Health"
$ws.Range("A2").Value = "₹1,514"
$ws.Range("B2").Value = "₹2,089"
$ws.Range("C2").Value = "₹2,156"
